# Added filtering options for the Component Analysis
#
# This trims the trailing "staircase" of comparison columns (G:K) on the
# early rows of the naive-error export so that each row only keeps the
# columns that fall inside the (now filtered) comparison window - mirroring
# the same diagonal cut-off pattern already used further down the sheet
# (e.g. rows 40-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2:K2").ClearContents()
$ws.Range("G3:K3").ClearContents()
$ws.Range("I4:K4").ClearContents()
$ws.Range("G5:K5").ClearContents()
$ws.Range("I6:K6").ClearContents()
$ws.Range("G7:K7").ClearContents()
$ws.Range("I8:K8").ClearContents()
$ws.Range("G9:K9").ClearContents()
$ws.Range("I10:K10").ClearContents()
$ws.Range("G11:K11").ClearContents()
$ws.Range("I12:K12").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("I14:K14").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("I16:K16").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("J18:K18").ClearContents()
$ws.Range("I19:K19").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("J22:K22").ClearContents()
$ws.Range("I23:K23").ClearContents()
$ws.Range("K25").ClearContents()
$ws.Range("J26:K26").ClearContents()
$ws.Range("I27:K27").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("J30:K30").ClearContents()
$ws.Range("I31:K31").ClearContents()
$ws.Range("K33").ClearContents()
$ws.Range("J34:K34").ClearContents()
$ws.Range("I35:K35").ClearContents()
$ws.Range("K37").ClearContents()
$ws.Range("J38:K38").ClearContents()
$ws.Range("I39:J39").ClearContents()
